# ------------------------------------------------------------------
# Add a new "2022-Q4" sheet (right after the "总计" summary sheet and
# before the existing "2022-Q3" sheet), record its two fund rows, and
# insert a matching summary row on the "总计" sheet.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------- 1. "总计" (summary) sheet: insert a new row for 2022-Q4 ----------
$summary = $wb.Worksheets.Item(1)

$summary.Rows(2).Insert()
$summary.Rows(2).ClearFormats()

# copy the style used by the other index cells in column A onto the new A2
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0.08

# re-number the index column for the rows that shifted down
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4

# ---------- 2. Create the "2022-Q4" detail sheet ----------
# Duplicate the existing "2022-Q3" sheet (position 2) so the new sheet
# starts out with identical layout/styling/fund code & name, then only
# update the figures that actually differ for 2022-Q4.
$q3 = $wb.Worksheets.Item(2)
$q3.Copy($q3)
$excel.CutCopyMode = 0

$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $q4.Range("D2") "5.40"
Set-TextValue $q4.Range("E2") "29.80"
Set-TextValue $q4.Range("F2") "1.49"
Set-TextValue $q4.Range("G2") "0.0805"
$q4.Range("H2").Value = 2

Set-TextValue $q4.Range("D3") "0.04"
Set-TextValue $q4.Range("E3") "29.80"
Set-TextValue $q4.Range("F3") "1.49"
Set-TextValue $q4.Range("G3") "0.0006"
$q4.Range("H3").Value = 2

Write-Host "done"
